$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$addr,
        [string]$val
    )
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "52.283.82"
$ws.Range("E2").Value = "  +0.85%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.911.79"
$ws.Range("E3").Value = "  +3.55%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.07%  "

# Row 5 - BNB
Set-TextValue "D5" "351.30"
$ws.Range("E5").Value = "  -1.56%  "

# Row 6 - Solana
Set-TextValue "D6" "112.41"
$ws.Range("E6").Value = "  +2.09%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +0.19%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.03%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.636"
$ws.Range("E9").Value = "  +0.29%  "

# Row 10 - Avalanche
Set-TextValue "D10" "40.04"
$ws.Range("E10").Value = "  -0.26%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.0864"
$ws.Range("E11").Value = "  +3.22%  "

# Row 12 - TRON
Set-TextValue "D12" "0.135"
$ws.Range("E12").Value = "  +0.04%  "

# Row 13 - Chainlink
Set-TextValue "D13" "20.00"
$ws.Range("E13").Value = "  -0.28%  "

# Row 14 - Polkadot
Set-TextValue "D14" "7.79"
$ws.Range("E14").Value = "  -0.13%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "3.373.70"
$ws.Range("E15").Value = "  +3.58%  "

# Row 16 - Polygon
$ws.Range("E16").Value = "  +6.87%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.923.81"
$ws.Range("E17").Value = "  +3.66%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "52.334.41"
$ws.Range("E18").Value = "  +0.94%  "

# Row 19 - was Uniswap, now ImmutableX
$ws.Range("B19").Value = "ImmutableX"
$ws.Range("C19").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D19" "3.32"
$ws.Range("E19").Value = "  +4.06%  "

# Row 20 - was ImmutableX, now Uniswap
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D20" "7.64"
$ws.Range("E20").Value = "  -0.42%  "

# Row 21 - InternetComputer(DFINITY)
Set-TextValue "D21" "14.21"
$ws.Range("E21").Value = "  +3.91%  "

# Row 22 - ShibaInu
$ws.Range("E22").Value = "  +0.17%  "

# Row 23 - Litecoin
Set-TextValue "D23" "71.02"
$ws.Range("E23").Value = "  +0.81%  "

# Row 24 - BitcoinCash
Set-TextValue "D24" "270.89"
$ws.Range("E24").Value = "  +0.88%  "

# Row 25
Set-TextValue "D25" "2.79"
$ws.Range("E25").Value = "  +1.16%  "

# Row 26
Set-TextValue "D26" "26.78"
$ws.Range("E26").Value = "  +2.31%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  -0.12%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  +1.69%  "

# Row 29 - Cosmos
Set-TextValue "D29" "10.63"
$ws.Range("E29").Value = "  +1.88%  "

# Row 30 - InjectiveProtocol
Set-TextValue "D30" "37.68"
$ws.Range("E30").Value = "  -0.72%  "

# Row 31 - was Toncoin, now RenderToken
$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D31" "6.25"
$ws.Range("E31").Value = "  +10.44%  "

# Row 32 - was Filecoin, now Toncoin
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D32" "2.25"
$ws.Range("E32").Value = "  +0.52%  "

# Row 33 - was RenderToken, now Filecoin
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D33" "6.48"
$ws.Range("E33").Value = "  +5.08%  "

# Row 34 - Hedera
Set-TextValue "D34" "0.0968"
$ws.Range("E34").Value = "  +11.24%  "

# Row 35 - OKB
$ws.Range("E35").Value = "  +1.99%  "

# Row 36 - VeChain
Set-TextValue "D36" "0.0451"
$ws.Range("E36").Value = "  +1.11%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  -0.11%  "

# Row 38 - LidoDAOToken
Set-TextValue "D38" "3.31"
$ws.Range("E38").Value = "  +5.38%  "

# Row 39 - was Celestia, now Stacks
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D39" "2.91"
$ws.Range("E39").Value = "  +16.55%  "

# Row 40 - was ARBITRUM, now Celestia
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D40" "18.84"
$ws.Range("E40").Value = "  +0.10%  "

# Row 41 - was Stacks, now ARBITRUM
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D41" "2.07"
$ws.Range("E41").Value = "  +2.99%  "

# Row 42 - Stellar
$ws.Range("E42").Value = "  +1.73%  "

# Row 43 - EnergySwap
Set-TextValue "D43" "23.61"
$ws.Range("E43").Value = "  +7.78%  "

# Row 44 - was ApeXProtocol, now Monero
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D44" "121.76"
$ws.Range("E44").Value = "  +1.74%  "

# Row 45 - was Monero, now ApeXProtocol
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D45" "2.63"
$ws.Range("E45").Value = "  +5.99%  "

# Row 46 - WEMIXToken
$ws.Range("E46").Value = "  -0.23%  "

# Row 47 - NEARProtocol
$ws.Range("E47").Value = "  +5.11%  "

# Row 48 - Maker
Set-TextValue "D48" "2.196.05"
$ws.Range("E48").Value = "  +4.24%  "

# Row 49 - TheGraph
Set-TextValue "D49" "0.264"
$ws.Range("E49").Value = "  +24.07%  "

# Row 50 - BEAM
Set-TextValue "D50" "0.0335"
$ws.Range("E50").Value = "  +10.18%  "

# Row 51 - SEI
Set-TextValue "D51" "0.961"
$ws.Range("E51").Value = "  +3.53%  "
